# Fill in the rule-index column (A3:A11) of the ServiceDetermination
# decision table with sequential text values "1".."9" (one per DMN rule
# row), matching the text-valued convention used by every other cell in
# this table (e.g. B3="-", M3="456" are all stored as text, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ServiceDetermination")

for ($row = 3; $row -le 11; $row++) {
    $ruleNumber = $row - 2
    $target = $ws.Cells.Item($row, 1)      # column A
    $formatSource = $ws.Cells.Item($row, 2) # column B, same row/style, already text

    # Force the cell to store the number as text (like the rest of the
    # table) instead of letting Excel auto-detect it as a numeric value.
    $target.NumberFormat = "@"
    $target.Value = [string]$ruleNumber

    # NumberFormat above stamps a brand-new style record, which would
    # bump A-column cells off the shared "s=3" style used throughout the
    # table. Restore the original formatting (border/font/number format)
    # by copying it from the neighboring same-row, same-style cell, while
    # leaving the freshly-written text value untouched.
    $formatSource.Copy()
    $target.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0
